# Insert a new data row at spreadsheet row 78 (pushes existing rows 78..158
# down to 79..159) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(78).Insert()

$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 44546
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 100112017
$ws.Range("G78").Value = "Apio"
$ws.Range("H78").Value = "Americana (o)"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 80
$ws.Range("K78").Value = 8000
$ws.Range("L78").Value = 8500
$ws.Range("M78").Value = 8250
$ws.Range("N78").Value = '$/docena de matas'
$ws.Range("O78").Value = "Provincia del Elquí"
$ws.Range("P78").Value = 1375
$ws.Range("Q78").Value = 6
$ws.Range("R78").Value = "Hortaliza"

Write-Output "done"
